$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 114.8270160096505
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 786.6765073101459

# Row 3 updated values
$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 4.429675500412797
